$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure the date/time text columns (E, F) are forced to plain text so that
# values such as "12/2/2026" or "5:20:17 pm" are NOT auto-converted into Excel
# date/time serial numbers - they must stay literal strings, exactly like the
# rest of the (already text-formatted) B/C/D columns.
$ws.Range("E2:F7").NumberFormat = "@"

# row, ID, Name, Email, Role, SignupDate, SignupTime
$rows = @(
    @(2, 4,  "ss",       "ss@gmail.com",      "student", "12/2/2026", "5:20:17 pm"),
    @(3, 3,  "vv",       "vv@gmail.com",      "student", "12/2/2026", "5:17:35 pm"),
    @(4, 2,  "vasu",     "vasu@gmail.com",    "admin",   "12/2/2026", "5:17:00 pm"),
    @(5, 9,  "vasu",     "vasuu@gmail.com",   "student", "17/2/2026", "5:54:19 pm"),
    @(6, 10, "vd",       "vd@gmail.com",      "student", "18/2/2026", "10:29:25 am"),
    @(7, 11, "John Doe", "john@example.com",  "student", "18/2/2026", "10:32:22 am")
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
    $ws.Cells.Item($rowNum, 5).Value = $r[5]
    $ws.Cells.Item($rowNum, 6).Value = $r[6]
}
